$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Result" invoice table lost an OCR'd row (invoice 284213) and the
# row that used to sit under it (invoice 284228) needs to disappear, with
# every row below it sliding up to close the gap.
#
# Insert a fresh row at row 2 (this pushes the existing rows 2-4 down to
# rows 3-5, keeping all of their original cell values/formatting intact),
# fill the new row 2 with the missing OCR result, and then delete what is
# now row 4 (the old "284228" row) so the table closes back up to 4 rows -
# leaving the old row 2 data at row 3 and the untouched old row 4 data
# (284232) back at row 4.
$ws.Rows(2).Insert()
$ws.Rows(4).Delete()

# New invoice row pulled in from OCR. Force these as text (matching how
# the rest of the "numeric-looking" invoice numbers / dates / amounts in
# this table are stored) so Excel doesn't silently convert them to a
# date serial / number.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"

$ws.Range("A2").Value = "284213"
$ws.Range("B2").Value = "2019-06-03"
$ws.Range("C2").Value = "Aenean LLC"
$ws.Range("D2").Value = "9778.40"
